$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.409.25"
$ws.Range("E2").Value = "  +2.19%  "

$ws.Range("D3").Value = "3.417.83"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.64"
$ws.Range("E5").Value = "  +1.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.97"
$ws.Range("E6").Value = "  +4.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +1.37%  "

$ws.Range("D8").Value = "3.409.01"
$ws.Range("E8").Value = "  +0.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.641"
$ws.Range("E11").Value = "  +1.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.79"
$ws.Range("E12").Value = "  +1.47%  "

$ws.Range("E13").Value = "  +0.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.34"
$ws.Range("E14").Value = "  +2.60%  "

$ws.Range("D15").Value = "3.965.06"
$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.38"
$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("D17").Value = "3.418.95"
$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("D19").Value = "66.336.96"
$ws.Range("E19").Value = "  +2.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.01"
$ws.Range("E20").Value = "  +1.55%  "

$ws.Range("E21").Value = "  +1.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "470.73"
$ws.Range("E22").Value = "  +0.55%  "

$ws.Range("E23").Value = "  +2.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.69"
$ws.Range("E24").Value = "  +8.25%  "

$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.57"
$ws.Range("E26").Value = "  +3.19%  "

$ws.Range("E27").Value = "  +0.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.88"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("E29").Value = "  +1.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.42"
$ws.Range("E30").Value = "  +2.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.95"
$ws.Range("E31").Value = "  +3.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.59"
$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "586.70"
$ws.Range("E33").Value = "  +2.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "62.58"
$ws.Range("E34").Value = "  +1.94%  "

$ws.Range("E35").Value = "  +1.23%  "

$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("E37").Value = "  +4.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.61"
$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.45"
$ws.Range("E39").Value = "  +2.44%  "

$ws.Range("E40").Value = "  +4.34%  "

$ws.Range("D41").Value = "0.0₃0759"
$ws.Range("E41").Value = "  +1.93%  "

$ws.Range("D42").Value = "3.133.36"
$ws.Range("E42").Value = "  +1.31%  "

$ws.Range("E43").Value = "  +2.37%  "

$ws.Range("E44").Value = "  +2.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.54"
$ws.Range("E45").Value = "  +3.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.80"
$ws.Range("E46").Value = "  +19.14%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.135"
$ws.Range("E47").Value = "  -0.11%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.20"
$ws.Range("E48").Value = "  +1.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.62"
$ws.Range("E50").Value = "  +2.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.61"
$ws.Range("E51").Value = "  +4.06%  "
